$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3874.25
$ws.Range("I74").Value = 3832.3333
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 3832.3333
$ws.Range("L74").Value = 4000
$ws.Range("M74").Value = -2896.3333
$ws.Range("N74").Value = -5872

$ws.Range("H77").Value = 3874.25
$ws.Range("I77").Value = 3832.3333
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 19161.6665
$ws.Range("L77").Value = 20000
$ws.Range("M77").Value = -14481.6665
$ws.Range("N77").Value = -29360

$ws.Range("H82").Value = 807
$ws.Range("I82").Value = 568.4
$ws.Range("K82").Value = 1705.2
$ws.Range("M82").Value = -1299.2

$ws.Range("H85").Value = 807
$ws.Range("I85").Value = 568.4
$ws.Range("K85").Value = 1705.2
$ws.Range("M85").Value = -301.1999999999998

$ws.Range("H141").Value = 9019.286
$ws.Range("I141").Value = 3275.25
$ws.Range("J141").Value = 11316.9
$ws.Range("K141").Value = 9825.75
$ws.Range("L141").Value = 33950.7
$ws.Range("M141").Value = -4645.75
$ws.Range("N141").Value = -44310.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 47780
$ws.Range("J76").Value = 47780
$ws.Range("L76").Value = 47780
$ws.Range("N76").Value = -48456

$ws.Range("H79").Value = 47780
$ws.Range("J79").Value = 47780
$ws.Range("L79").Value = 47780
$ws.Range("N79").Value = -50120

$ws.Range("H132").Value = 1435919
$ws.Range("I132").Value = 976.4722
$ws.Range("K132").Value = 2929.4166
$ws.Range("M132").Value = -399.4166

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()

$ws.Range("H35").Value = 24000
$ws.Range("J35").Value = 24000
$ws.Range("L35").Value = 24000
$ws.Range("N35").Value = -24620

$ws.Range("H76").Value = 14500
$ws.Range("I76").Value = 4000
$ws.Range("K76").Value = 4000
$ws.Range("M76").Value = -3685

$ws.Range("H79").Value = 14500
$ws.Range("I79").Value = 4000
$ws.Range("K79").Value = 4000
$ws.Range("M79").Value = -2908

$ws.Range("H94").Value = 1137.1428
$ws.Range("I94").Value = 1660
$ws.Range("J94").Value = 440
$ws.Range("K94").Value = 1660
$ws.Range("L94").Value = 440
$ws.Range("M94").Value = -1209
$ws.Range("N94").Value = -1342

$ws.Range("H138").Value = 56629.168
$ws.Range("J138").Value = 56629.168
$ws.Range("L138").Value = 56629.168
$ws.Range("N138").Value = -66909.16800000001

$ws.Range("H139").Value = 46114.082
$ws.Range("I139").Value = 49980
$ws.Range("J139").Value = 45762.637
$ws.Range("K139").Value = 49980
$ws.Range("L139").Value = 45762.637
$ws.Range("M139").Value = -44840
$ws.Range("N139").Value = -56042.637

$ws.Range("H140").Value = 37464.668
$ws.Range("J140").Value = 37464.668
$ws.Range("L140").Value = 37464.668
$ws.Range("N140").Value = -47824.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5557246.5
$ws.Range("I31").Value = 9260044
$ws.Range("J31").Value = 3050
$ws.Range("K31").Value = 9260044
$ws.Range("L31").Value = 3050
$ws.Range("M31").Value = -9259749
$ws.Range("N31").Value = -3640

$ws.Range("H34").Value = 5557246.5
$ws.Range("I34").Value = 9260044
$ws.Range("J34").Value = 3050
$ws.Range("K34").Value = 9260044
$ws.Range("L34").Value = 3050
$ws.Range("M34").Value = -9259842
$ws.Range("N34").Value = -3454

$ws.Range("H56").Value = 25000
$ws.Range("J56").Value = 25000
$ws.Range("L56").Value = 25000
$ws.Range("N56").Value = -26690

$ws.Range("H103").Value = 6500
$ws.Range("I103").Value = 6500
$ws.Range("K103").Value = 6500
$ws.Range("M103").Value = -5328

$ws.Range("H139").Value = 49587.9
$ws.Range("J139").Value = 49587.9
$ws.Range("L139").Value = 49587.9
$ws.Range("N139").Value = -59867.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 170
$ws.Range("I15").Value = 170
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 510
$ws.Range("L15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = -370

$ws.Range("H25").Value = 90
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 90
$ws.Range("K25").Value = 0
$ws.Range("L25").ClearContents()
$ws.Range("M25").Value = 270
$ws.Range("N25").Value = -608

$ws.Range("H30").Value = 90
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 90
$ws.Range("K30").Value = 0
$ws.Range("L30").ClearContents()
$ws.Range("M30").Value = 270
$ws.Range("N30").Value = -474

$ws.Range("H137").Value = 14286979
$ws.Range("I137").Value = 38462280
$ws.Range("J137").Value = 1574.2273
$ws.Range("K137").Value = 115386840
$ws.Range("L137").Value = 4722.6819
$ws.Range("M137").Value = -115381740
$ws.Range("N137").Value = -14922.6819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2243
$ws.Range("I102").Value = 1562.3
$ws.Range("K102").Value = 1562.3
$ws.Range("M102").Value = 59.70000000000005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H8").Value = 32000
$ws.Range("J8").Value = 32000
$ws.Range("L8").Value = 32000
$ws.Range("N8").Value = -32280

$ws.Range("H80").Value = 58835
$ws.Range("J80").Value = 58835
$ws.Range("L80").Value = 58835
$ws.Range("N80").Value = -61081

$ws.Range("H83").Value = 58835
$ws.Range("J83").Value = 58835
$ws.Range("L83").Value = 176505
$ws.Range("N83").Value = -187737

$ws.Range("H93").Value = 1439.125
$ws.Range("I93").Value = 1191.7778
$ws.Range("J93").Value = 1757.1428
$ws.Range("K93").Value = 1191.7778
$ws.Range("L93").Value = 1757.1428
$ws.Range("M93").Value = 56.22219999999993
$ws.Range("N93").Value = -4253.1428

$ws.Range("H100").Value = 1837.5186
$ws.Range("I100").Value = 3196.6
$ws.Range("J100").Value = 1528.6364
$ws.Range("K100").Value = 3196.6
$ws.Range("L100").Value = 1528.6364
$ws.Range("M100").Value = -2655.6
$ws.Range("N100").Value = -2610.6364

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1127.5714
$ws.Range("J81").Value = 1363.5
$ws.Range("L81").Value = 2727
$ws.Range("N81").Value = -4849

$ws.Range("H84").Value = 1127.5714
$ws.Range("J84").Value = 1363.5
$ws.Range("L84").Value = 13635
$ws.Range("N84").Value = -24243
